$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new 2019 data row (year, Sept, Oct, Nov, Dec, Total) after the
# existing last row (2018, row 53).
$ws.Range("A54").Value = 2019
$ws.Range("B54").Value = 0
$ws.Range("C54").Value = 0
$ws.Range("D54").Value = 0
$ws.Range("E54").Value = 0
$ws.Range("F54").Value = 0

# Match the selection left behind after entering the new row of data.
$ws.Range("A55").Select()
